# Update capital structure database values for Panama / Financial Svcs.
# (Non-bank & Insurance) rows (rows 2 and 3 share identical figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3)

foreach ($r in $rows) {
    $ws.Range("D$r").Value  = -0.0795
    $ws.Range("E$r").Value  = -0.09039999999999999

    $ws.Range("I$r").Value  = 0
    $ws.Range("J$r").Value  = 0
    $ws.Range("K$r").Value  = 69.90000000000001
    $ws.Range("L$r").Value  = 0.6424632352941178
    $ws.Range("M$r").Value  = 49.8
    $ws.Range("N$r").Value  = 0.07929936305732484
    $ws.Range("O$r").Value  = 0.7124463519313303
    $ws.Range("P$r").Value  = 49.8
    $ws.Range("Q$r").Value  = 0.07929936305732484
    $ws.Range("R$r").Value  = 0.7124463519313303
    $ws.Range("S$r").Value  = 0
    $ws.Range("T$r").Value  = 0
    $ws.Range("U$r").Value  = 1357.3
    $ws.Range("V$r").Value  = 2.161305732484077
    $ws.Range("W$r").Value  = 0.06926964621940343
    $ws.Range("X$r").Value  = 0.04553149513255872
    $ws.Range("Y$r").Value  = 0.02373815108684471
    $ws.Range("Z$r").Value  = 0.03969064643221946
    $ws.Range("AA$r").Value = 0
    $ws.Range("AB$r").Value = 0.03385338033707589
    $ws.Range("AC$r").Value = -0.03385338033707589
    $ws.Range("AD$r").Value = 2110.9
    $ws.Range("AE$r").Value = 0
    $ws.Range("AF$r").Value = 2110.9
    $ws.Range("AG$r").Value = 753.6000000000001
    $ws.Range("AH$r").Value = 0.7707108693271021
    $ws.Range("AI$r").Value = 0.672989861633616
    $ws.Range("AJ$r").Value = 0.5454545454545455
    $ws.Range("AK$r").Value = 0.4235373461473614

    # debt_ebitda / net_debt_ebitda are no longer populated for this row.
    $ws.Range("AN$r").ClearContents()
    $ws.Range("AP$r").ClearContents()
}
